# ---------------------------------------------------------------------------
# 1) Slide 16 table: change the applied table style GUID.
# ---------------------------------------------------------------------------
$p = $ppt.ActivePresentation

$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{763EC968-6CFA-434D-AF91-1CF6D7E1EF36}")
    }
}

# ---------------------------------------------------------------------------
# 2) Theme colour swap: the deck's theme ("Integral") and the notes-master's
#    theme ("Office Theme") had their colour palettes swapped. The slide
#    master's live theme colour scheme is reachable through COM -- set each
#    of its 12 slots to the Office-Theme palette values (encoded the way
#    VBA's RGB()/RGBColor.RGB represent colour: R + G*256 + B*65536, i.e.
#    hex literal 0xBBGGRR).
# ---------------------------------------------------------------------------
$cs = $p.SlideMaster.Theme.ThemeColorScheme

$cs.Item(1).RGB  = 0x000000   # dk1      000000
$cs.Item(2).RGB  = 0xFFFFFF   # lt1      FFFFFF
$cs.Item(3).RGB  = 0x6A5444   # dk2      44546A
$cs.Item(4).RGB  = 0xE6E6E7   # lt2      E7E6E6
$cs.Item(5).RGB  = 0xD59B5B   # accent1  5B9BD5
$cs.Item(6).RGB  = 0x317DED   # accent2  ED7D31
$cs.Item(7).RGB  = 0xA5A5A5   # accent3  A5A5A5
$cs.Item(8).RGB  = 0x00C0FF   # accent4  FFC000
$cs.Item(9).RGB  = 0xC47244   # accent5  4472C4
$cs.Item(10).RGB = 0x47AD70   # accent6  70AD47
$cs.Item(11).RGB = 0xC16305   # hlink    0563C1
$cs.Item(12).RGB = 0x724F95   # folHlink 954F72
